$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Unveiling Quantum Secrets: A Journey into the Unseen" "Biology: Exploring the Essence of Life"

# Author title "Prof" -> "Mrs"
Replace-Text "Prof" "Mrs"

# Author name: " Tamara Morgan" -> " Emily C" (rest appended as new runs below)
Replace-Text " Tamara Morgan" " Emily C"

# Email user/domain parts
Replace-Text "tamara" "emily"
Replace-Text "morgan01@ashevillecollege" "richards@biotech"
Replace-Text "net" "edu"

# Body paragraph 1
Replace-Text "In the vast tapestry of the universe lies a realm of secrets that has captivated humankind for centuries: the quantum world" "Biology, the study of life, unravels the intricacies of existence in all its forms, from minuscule microorganisms to towering canopies"
Replace-Text " An enigmatic realm where particles behave in ways that defy classical physics, quantum mechanics promises to revolutionize our understanding of the universe and open up new frontiers of technology" " It illuminates not just the fractional components of living systems, but unravels the interconnected narrative that orchestrates their harmony"
Replace-Text " Delving into this uncharted territory, scientists and philosophers embark on a journey to unlock the hidden mysteries and harness the untapped potential of the quantum realm" " Beyond the confines of textbooks and laboratory walls, biology unveils a universe of knowledge woven into the fabric of existence"

Replace-Text "Unveiling the secrets of the quantum world requires delving into the depths of theoretical physics, where equations dance in harmony to describe the behavior of subatomic particles" "Biology illuminates the foundation of life's complexity- cells, the fundamental building blocks of all living organisms"
Replace-Text " The enigma of quantum entanglement, where particles located miles apart can instantaneously influence each other, continues to perplex scientists and has far-reaching implications for communication and computation" " Within these microscopic realms, intercellular communication, complex biochemical reactions, and the symphony of genetic information orchestrate a symphony of life"
Replace-Text " As researchers navigate the intricate landscape of quantum mechanics, they are discovering phenomena such as superposition, where particles exist in multiple states simultaneously, shattering our conventional notions of reality" " The study of cells sparks curiosity into the very essence of living, unveiling the codes of heredity and the secrets to longevity and disease"

Replace-Text "The implications of quantum mechanics extend far beyond the realm of theory and into the realm of practical application" "Biology explores ecosystems, where organisms and environments intertwine"
Replace-Text " Quantum technologies hold the promise of revolutionizing diverse fields, ranging from secure communication to ultra-precise measurements" " The rich tapestry of life, teeming with biodiversity, exhibits exquisite patterns of interaction among its members"
Replace-Text " Quantum cryptography, for instance, harnesses the laws of quantum mechanics to encode information in ways that are inherently secure, offering unprecedented levels of protection against eavesdropping" " From the microscopic world of decomposition to the grand procession of migration, nature displays intricate collaboration and competition"
Replace-Text " Quantum computers, with their ability to perform calculations exponentially faster than conventional computers, promise to unlock new frontiers of innovation in fields such as drug discovery and materials science" " Understanding ecosystems empowers us to comprehend the intricate balance of nature and grapple with urgent global issues like climate change"

# Summary paragraph
Replace-Text "In the tapestry of the universe, the quantum world stands as a realm of exquisite mystery and profound implications" "Biology embarks on a journey to comprehend the intricacies of life, from cells to ecosystems"
Replace-Text " As scientists delve deeper into its secrets, they are uncovering insights that challenge our conventional understanding of reality and ushering in a new era of technological innovation" " It deciphers the language of genetic information, unravels the structure and function of organisms, and paints a vivid tapestry of living systems interacting with their environments"
Replace-Text " From the enigmatic realm of quantum entanglement to the promise of quantum technologies, the journey into the unseen realm of quantum mechanics is transforming our understanding of the universe and shaping the future of human ingenuity" " By unveiling the interconnectedness of life, biology equips us to appreciate the richness of our world and address global challenges"

# Split the author name run " Emily C" into " Emily C" + "." + " Richards" as separate runs
$d.Content.Find.Execute(" Emily C", $false, $true) | Out-Null
$found = $d.Content.Find.Found
$nameRange = $d.Content.Duplicate
$nameRange.Find.Execute(" Emily C", $false, $true) | Out-Null
$nameRange.Collapse(0)
$nameRange.InsertAfter(".")
$nameRange.Collapse(0)
$nameRange.InsertAfter(" Richards")

# Append a new empty paragraph at the end of the document body
$d.Paragraphs.Last.Range.InsertParagraphAfter()
